$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.750.55'
$ws.Range('E2').Value = '  -5.40%  '
$ws.Range('D3').Value = '3.299.56'
$ws.Range('E3').Value = '  -6.76%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = '556.98'
$ws.Range('E5').Value = '  -5.12%  '
$ws.Range('D6').Value = '180.40'
$ws.Range('E6').Value = '  -7.30%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '0.589'
$ws.Range('E8').Value = '  -3.28%  '
$ws.Range('D9').Value = '3.292.91'
$ws.Range('E9').Value = '  -6.58%  '
$ws.Range('E10').Value = '  -8.42%  '
$ws.Range('D11').Value = '0.587'
$ws.Range('E11').Value = '  -5.63%  '
$ws.Range('D12').Value = '47.68'
$ws.Range('E12').Value = '  -9.35%  '
$ws.Range('D13').Value = '0.0000265'
$ws.Range('E13').Value = '  -8.01%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = '8.56'
$ws.Range('E14').Value = '  -7.32%  '
$ws.Range('B15').Value = 'BitcoinCash'
$ws.Range('C15').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D15').Value = '630.02'
$ws.Range('E15').Value = '  -2.74%  '
$ws.Range('D16').Value = '3.825.55'
$ws.Range('E16').Value = '  -6.90%  '
$ws.Range('D17').Value = '18.04'
$ws.Range('E17').Value = '  -1.92%  '
$ws.Range('D18').Value = '65.709.56'
$ws.Range('E18').Value = '  -5.60%  '
$ws.Range('E19').Value = '  -4.07%  '
$ws.Range('D20').Value = '3.290.29'
$ws.Range('E20').Value = '  -7.05%  '
$ws.Range('D21').Value = '11.43'
$ws.Range('E21').Value = '  -9.06%  '
$ws.Range('D22').Value = '0.907'
$ws.Range('E22').Value = '  -5.47%  '
$ws.Range('D23').Value = '17.67'
$ws.Range('E23').Value = '  -2.17%  '
$ws.Range('D24').Value = '106.26'
$ws.Range('E24').Value = '  +3.84%  '
$ws.Range('D25').Value = '5.04'
$ws.Range('E25').Value = '  -8.25%  '
$ws.Range('D26').Value = '4.00'
$ws.Range('E26').Value = '  -8.51%  '
$ws.Range('D27').Value = '6.01'
$ws.Range('E27').Value = '  -0.44%  '
$ws.Range('E28').Value = '  -8.13%  '
$ws.Range('D29').Value = '9.51'
$ws.Range('E29').Value = '  -6.53%  '
$ws.Range('D30').Value = '8.78'
$ws.Range('E30').Value = '  -8.11%  '
$ws.Range('D31').Value = '30.66'
$ws.Range('E31').Value = '  -7.33%  '
$ws.Range('D32').Value = '4.02'
$ws.Range('E32').Value = '  -3.75%  '
$ws.Range('D33').Value = '6.30'
$ws.Range('E33').Value = '  -7.00%  '
$ws.Range('D34').Value = '11.07'
$ws.Range('E34').Value = '  -5.69%  '
$ws.Range('B35').Value = 'Bittensor'
$ws.Range('C35').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D35').Value = '547.71'
$ws.Range('E35').Value = '  +7.06%  '
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').Value = '0.105'
$ws.Range('E36').Value = '  -4.66%  '
$ws.Range('D37').Value = '3.738.54'
$ws.Range('E37').Value = '  +0.31%  '
$ws.Range('D38').Value = '0.999'
$ws.Range('E38').Value = '  +0.04%  '
$ws.Range('D39').Value = '56.81'
$ws.Range('E39').Value = '  -7.77%  '
$ws.Range('D40').Value = '3.47'
$ws.Range('E40').Value = '  -4.35%  '
$ws.Range('B41').Value = 'Fetch.AI'
$ws.Range('C41').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D41').Value = '2.76'
$ws.Range('E41').Value = '  -6.53%  '
$ws.Range('B42').Value = 'PEPE'
$ws.Range('C42').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D42').Value = '0.0₃0711'
$ws.Range('E42').Value = '  -11.84%  '
$ws.Range('B43').Value = 'CoreDAO'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D43').Value = '3.38'
$ws.Range('E43').Value = '  +24.87%  '
$ws.Range('E44').Value = '  -6.30%  '
$ws.Range('D45').Value = '0.342'
$ws.Range('E45').Value = '  -7.63%  '
$ws.Range('D46').Value = '32.02'
$ws.Range('E46').Value = '  -7.93%  '
$ws.Range('D47').Value = '0.0413'
$ws.Range('E47').Value = '  -7.64%  '
$ws.Range('D48').Value = '3.22'
$ws.Range('E48').Value = '  -6.12%  '
$ws.Range('B49').Value = 'ThetaToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D49').Value = '2.62'
$ws.Range('E49').Value = '  -8.51%  '
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').Value = '0.129'
$ws.Range('E50').Value = '  -4.91%  '
